$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("UWKT3")

# Michal Monselise withdrew from the team but the roster keeps her row;
# a blank row is inserted above the roster table (row 5), shifting every
# member row down by one.
$ws.Rows("5:5").Insert()

# The hyperlink on Gregory Hogue's e-mail address moved down with his row
# (old H7 -> new H8). Re-attach it so the link still points at the right cell.
$ws.Range("H7").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("H8"), "mailto:gphogue@gmail.com")

# Werner Colangelo is added as a new team member, working alone, on a new
# row 17 (row 15 - the "YB" legend note - is now row 16). He only has a
# Full/First/Last name and an e-mail address; no Netid/Team/Online/Initials.
$ws.Range("A17").Value = "Werner Colangelo"
$ws.Range("B17").Value = "Werner"
$ws.Range("C17").Value = "Colangelo"
$ws.Range("H17").Value = "wernercolangelo@gmail.com"

# Match the author's final selection as recorded in the saved file.
$ws.Range("G24").Select()
